$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5299797058105469
$ws.Range("E2").Value = 1019.70200244529
$ws.Range("F2").Value = 0.05574637620534637
$ws.Range("G2").Value = 0.04437923226334865
$ws.Range("H2").Value = 0.03772822871291033
$ws.Range("I2").Value = 0.03109668840907373
$ws.Range("J2").Value = 0.02985795354131299
$ws.Range("K2").Value = 0.02782740690708877
$ws.Range("L2").Value = 0.02527091369511278
$ws.Range("M2").Value = 0.02450538904166488
$ws.Range("N2").Value = 0.02387691084537536
$ws.Range("O2").Value = 0.02297260866437747
$ws.Range("P2").Value = 0.02229808911495085
$ws.Range("Q2").Value = 0.02179792141590256
$ws.Range("R2").Value = 0.02144499790375659
$ws.Range("S2").Value = 0.02099916720422563
$ws.Range("T2").Value = 0.02042159354486529
$ws.Range("U2").Value = 0.02021302872845707
$ws.Range("V2").Value = 0.0201297555755928
$ws.Range("W2").Value = 0.02007244423370479
$ws.Range("X2").Value = 0.01998643468226164
$ws.Range("Y2").Value = 0.01987723201647739

$ws.Range("C3").Value = 0.5624589920043945
$ws.Range("E3").Value = 1035.003525387427
$ws.Range("F3").Value = 0.05592175031906621
$ws.Range("G3").Value = 0.04503869969666358
$ws.Range("H3").Value = 0.03872672915448062
$ws.Range("I3").Value = 0.03305653373628736
$ws.Range("J3").Value = 0.02972066225805825
$ws.Range("K3").Value = 0.02720865533002883
$ws.Range("L3").Value = 0.0250606663408576
$ws.Range("M3").Value = 0.02398243640209744
$ws.Range("N3").Value = 0.02336551052206967
$ws.Range("O3").Value = 0.02228903468281678
$ws.Range("P3").Value = 0.02187927634280718
$ws.Range("Q3").Value = 0.02154487650747213
$ws.Range("R3").Value = 0.02137877325810643
$ws.Range("S3").Value = 0.0208338978801405
$ws.Range("T3").Value = 0.0207087325597422
$ws.Range("U3").Value = 0.02059503258302952
$ws.Range("V3").Value = 0.02049387303284274
$ws.Range("W3").Value = 0.02037805596261094
$ws.Range("X3").Value = 0.02027740795994096
$ws.Range("Y3").Value = 0.02017550731749369

$ws.Range("C4").Value = 0.6526656150817871
$ws.Range("E4").Value = 1033.424683337909
$ws.Range("F4").Value = 0.05564009632728495
$ws.Range("G4").Value = 0.04354744946236904
$ws.Range("H4").Value = 0.03687403323858946
$ws.Range("I4").Value = 0.03258620175287177
$ws.Range("J4").Value = 0.02990718936213414
$ws.Range("K4").Value = 0.02752545336459601
$ws.Range("L4").Value = 0.02597640403475083
$ws.Range("M4").Value = 0.02525711690252885
$ws.Range("N4").Value = 0.02402188470921438
$ws.Range("O4").Value = 0.02281350020176161
$ws.Range("P4").Value = 0.02221199659611388
$ws.Range("Q4").Value = 0.02151318722595409
$ws.Range("R4").Value = 0.02120802539657932
$ws.Range("S4").Value = 0.02111358722965319
$ws.Range("T4").Value = 0.02078530232559442
$ws.Range("U4").Value = 0.02068776288762688
$ws.Range("V4").Value = 0.02036806864636418
$ws.Range("W4").Value = 0.02029491206856218
$ws.Range("X4").Value = 0.02019228356120901
$ws.Range("Y4").Value = 0.02014473066935495

$ws.Range("C5").Value = 0.5469143390655518
$ws.Range("E5").Value = 1063.546788421147
$ws.Range("F5").Value = 0.0568949876258323
$ws.Range("G5").Value = 0.04031775537175073
$ws.Range("H5").Value = 0.0373899473416082
$ws.Range("I5").Value = 0.03056792504716855
$ws.Range("J5").Value = 0.02817153958539872
$ws.Range("K5").Value = 0.02614937031773353
$ws.Range("L5").Value = 0.02460222578907969
$ws.Range("M5").Value = 0.0235589348662601
$ws.Range("N5").Value = 0.02305193100029102
$ws.Range("O5").Value = 0.02256772245392924
$ws.Range("P5").Value = 0.02200739658919712
$ws.Range("Q5").Value = 0.02171231597405218
$ws.Range("R5").Value = 0.02162131763673817
$ws.Range("S5").Value = 0.02142743126237694
$ws.Range("T5").Value = 0.0213466598890132
$ws.Range("U5").Value = 0.02109019667127613
$ws.Range("V5").Value = 0.0209871096689214
$ws.Range("W5").Value = 0.02088912152275559
$ws.Range("X5").Value = 0.02078077828494377
$ws.Range("Y5").Value = 0.02073190620703989

$ws.Range("C6").Value = 0.5312392711639404
$ws.Range("E6").Value = 1036.118830085525
$ws.Range("F6").Value = 0.05690075257779623
$ws.Range("G6").Value = 0.04412507034628728
$ws.Range("H6").Value = 0.03935646995284391
$ws.Range("I6").Value = 0.03292328026609049
$ws.Range("J6").Value = 0.03128124430848177
$ws.Range("K6").Value = 0.02900267881554905
$ws.Range("L6").Value = 0.02724013937342639
$ws.Range("M6").Value = 0.02546478869605168
$ws.Range("N6").Value = 0.02388498847207237
$ws.Range("O6").Value = 0.02248471881598893
$ws.Range("P6").Value = 0.02248471881598893
$ws.Range("Q6").Value = 0.02204649361127386
$ws.Range("R6").Value = 0.02165274924881268
$ws.Range("S6").Value = 0.02135492653951438
$ws.Range("T6").Value = 0.02091017814813053
$ws.Range("U6").Value = 0.02082574061839027
$ws.Range("V6").Value = 0.0205261673082386
$ws.Range("W6").Value = 0.02040471947780422
$ws.Range("X6").Value = 0.02029660682536025
$ws.Range("Y6").Value = 0.02019724814981529

$ws.Range("C7").Value = 0.5312492847442627
$ws.Range("E7").Value = 1035.072143769597
$ws.Range("F7").Value = 0.05533056042549669
$ws.Range("G7").Value = 0.04206417198151873
$ws.Range("H7").Value = 0.03544704245960476
$ws.Range("I7").Value = 0.0323793742834693
$ws.Range("J7").Value = 0.02888594684055329
$ws.Range("K7").Value = 0.02667678544153004
$ws.Range("L7").Value = 0.02484316089988795
$ws.Range("M7").Value = 0.02370465432152377
$ws.Range("N7").Value = 0.02287665274178744
$ws.Range("O7").Value = 0.02282403770281564
$ws.Range("P7").Value = 0.02228142940749529
$ws.Range("Q7").Value = 0.02185753351322837
$ws.Range("R7").Value = 0.02149944621668582
$ws.Range("S7").Value = 0.02110417626155506
$ws.Range("T7").Value = 0.02085237754081659
$ws.Range("U7").Value = 0.02069828468280987
$ws.Range("V7").Value = 0.02035823824621688
$ws.Range("W7").Value = 0.02029059350227894
$ws.Range("X7").Value = 0.02024901104846582
$ws.Range("Y7").Value = 0.02017684490778942

$ws.Range("C8").Value = 0.5312643051147461
$ws.Range("E8").Value = 1011.374670926971
$ws.Range("F8").Value = 0.0570572511652736
$ws.Range("G8").Value = 0.0439577586025639
$ws.Range("H8").Value = 0.03634328275752868
$ws.Range("I8").Value = 0.03132400532627549
$ws.Range("J8").Value = 0.02878502967848924
$ws.Range("K8").Value = 0.0268742208372783
$ws.Range("L8").Value = 0.02492156938760278
$ws.Range("M8").Value = 0.02307432123487167
$ws.Range("N8").Value = 0.02233282366694014
$ws.Range("O8").Value = 0.02217077735999274
$ws.Range("P8").Value = 0.02129444374671664
$ws.Range("Q8").Value = 0.02091137235778121
$ws.Range("R8").Value = 0.02033841442701992
$ws.Range("S8").Value = 0.02012938602753276
$ws.Range("T8").Value = 0.02012938602753276
$ws.Range("U8").Value = 0.01998249283353178
$ws.Range("V8").Value = 0.01989764021239524
$ws.Range("W8").Value = 0.0198086934038748
$ws.Range("X8").Value = 0.0197530727407901
$ws.Range("Y8").Value = 0.01971490586602282

$ws.Range("C9").Value = 0.546860933303833
$ws.Range("E9").Value = 1084.611399967262
$ws.Range("F9").Value = 0.05750392142240357
$ws.Range("G9").Value = 0.04281926014358252
$ws.Range("H9").Value = 0.03787745822595858
$ws.Range("I9").Value = 0.03408283620459283
$ws.Range("J9").Value = 0.03158727529567837
$ws.Range("K9").Value = 0.02882406738522839
$ws.Range("L9").Value = 0.02807490901861303
$ws.Range("M9").Value = 0.02587582508782843
$ws.Range("N9").Value = 0.02424951095364666
$ws.Range("O9").Value = 0.02379350527414254
$ws.Range("P9").Value = 0.02341354899532801
$ws.Range("Q9").Value = 0.02277730864812228
$ws.Range("R9").Value = 0.02244332928543035
$ws.Range("S9").Value = 0.02196236065885957
$ws.Range("T9").Value = 0.02176112971780946
$ws.Range("U9").Value = 0.02155536649327599
$ws.Range("V9").Value = 0.02147914426930909
$ws.Range("W9").Value = 0.02138631101929624
$ws.Range("X9").Value = 0.0212481931945412
$ws.Range("Y9").Value = 0.02114252241651582

$ws.Range("C10").Value = 0.5468754768371582
$ws.Range("E10").Value = 1040.785534694567
$ws.Range("F10").Value = 0.05761312742640419
$ws.Range("G10").Value = 0.04500057984911927
$ws.Range("H10").Value = 0.03816266930646862
$ws.Range("I10").Value = 0.03285991188474895
$ws.Range("J10").Value = 0.03034908270770848
$ws.Range("K10").Value = 0.02829960182491375
$ws.Range("L10").Value = 0.02568042270553055
$ws.Range("M10").Value = 0.02514750068825463
$ws.Range("N10").Value = 0.02322430947923576
$ws.Range("O10").Value = 0.02260217474920155
$ws.Range("P10").Value = 0.02182141918387152
$ws.Range("Q10").Value = 0.02181488645377573
$ws.Range("R10").Value = 0.0214450712299814
$ws.Range("S10").Value = 0.02106938496136179
$ws.Range("T10").Value = 0.02087129863822416
$ws.Range("U10").Value = 0.02078797996051353
$ws.Range("V10").Value = 0.02065047816959502
$ws.Range("W10").Value = 0.02043844840119715
$ws.Range("X10").Value = 0.02039310163327743
$ws.Range("Y10").Value = 0.02028821705057635

$ws.Range("C11").Value = 0.5468752384185791
$ws.Range("E11").Value = 1032.943466937362
$ws.Range("F11").Value = 0.05740972819423515
$ws.Range("G11").Value = 0.04701174230212647
$ws.Range("H11").Value = 0.03887248331927829
$ws.Range("I11").Value = 0.03316300059999087
$ws.Range("J11").Value = 0.0295697074415427
$ws.Range("K11").Value = 0.02686255809006518
$ws.Range("L11").Value = 0.02621640932262759
$ws.Range("M11").Value = 0.02489370009727366
$ws.Range("N11").Value = 0.02362597960902373
$ws.Range("O11").Value = 0.02269667542996461
$ws.Range("P11").Value = 0.02213834163472707
$ws.Range("Q11").Value = 0.02159560652139055
$ws.Range("R11").Value = 0.02118279027441449
$ws.Range("S11").Value = 0.0208820513114055
$ws.Range("T11").Value = 0.0207010071610337
$ws.Range("U11").Value = 0.02055390253073068
$ws.Range("V11").Value = 0.0204095419016854
$ws.Range("W11").Value = 0.02026379334276661
$ws.Range("X11").Value = 0.02020648377732294
$ws.Range("Y11").Value = 0.02013535023269711
